# Add AUTH-02 test case row to the "Authorization Module" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Authorization Module")

# Widen column D slightly to fit the new, longer description text.
$ws.Columns.Item(4).ColumnWidth = 48.09

# New row 5 values (mirrors the layout of row 4). Set in the same order the
# new shared strings were authored so the shared-strings table lines up.
$ws.Range("C5").Value = "AUTH-02"
$ws.Range("E5").Value = "AUTH-TC-02"
$ws.Range("G5").Value = "AUTH-TS-02"
$ws.Range("D5").Value = "Verify redirect to login when accessing`n restricted page via direct URL"
$ws.Range("H5").Value = "Access restricted page via direct `nURL without login"
$ws.Range("J5").Value = "1. Copy restricted page URL `n2. Paste into browser `n3. Press Enter"
$ws.Range("L5").Value = "User redirected to login page"
$ws.Range("F5").Value = "Authorization Module"
$ws.Range("I5").Value = "User not logged in"
$ws.Range("K5").Value = "Direct restricted URL"
$ws.Range("M5").Value = "High"

# Copy the formatting of row 4 onto row 5 so borders/fills/alignment match.
$ws.Range("C4:M4").Copy()
$ws.Range("C5:M5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# D5 keeps D4's fill/border/alignment (from the paste above) but also wraps
# its longer text.
$ws.Range("D5").WrapText = $true

$ws.Rows.Item(5).RowHeight = 63

# Scroll the viewport right and leave the selection where the author's
# session ended up after reviewing the new row.
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("L6").Select()
